$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ind")

# The SMA indicator route now documents support for daily & minutely
# timeframes (not just H1), and the example "options" query value changes
# from 10 to 2.

# C11: example request URL options=10 -> options=2
$ws.Range("C11").Value = "http://localhost:8080/ind?type=sma&options=2&pair=BTC,USD&timeframe=H1"

# D9: "Currently supported: H1" -> expanded timeframe list
$ws.Range("D9").Value = "Currently supported: D30, D1, H1, H2, H4, M30, M15"

# C13: Response was "(working on it)", now points to the status doc.
# Also drop the one-off font override this cell had, reverting to the
# sheet's normal/default style.
$ws.Range("C13").Value = "( status in indicators.xlsx)"
$ws.Range("C13").Style = "Standard"

# Row 14 (which used to hold "( status in indicators.xlsx)" on its own line)
# is no longer needed now that C13 carries that text directly.
$ws.Range("C14").Delete()

# Update the active selection to reflect where the edit left off.
$ws.Range("D9").Select()
